# Update "想去人数" (F column) figures that changed between crawls.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1901
$ws1.Range("F4").Value  = 892
$ws1.Range("F5").Value  = 790
$ws1.Range("F6").Value  = 13396
$ws1.Range("F7").Value  = 13263
$ws1.Range("F19").Value = 403
$ws1.Range("F20").Value = 279
$ws1.Range("F21").Value = 294
$ws1.Range("F22").Value = 429
$ws1.Range("F24").Value = 22

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 35

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1901
$ws4.Range("F5").Value  = 892
$ws4.Range("F7").Value  = 790
$ws4.Range("F8").Value  = 13396
$ws4.Range("F9").Value  = 13263
$ws4.Range("F26").Value = 403
$ws4.Range("F27").Value = 279
$ws4.Range("F28").Value = 294
$ws4.Range("F29").Value = 429
$ws4.Range("F33").Value = 22
$ws4.Range("F34").Value = 35
